$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.1158
$ws.Range("B4").Value = -0.0233
$ws.Range("B5").Value = 0.1614
$ws.Range("B6").Value = -0.0606
$ws.Range("B7").Value = 0.066
$ws.Range("B8").Value = -0.8298
$ws.Range("B9").Value = 0.0209
$ws.Range("B10").Value = 0.0402
$ws.Range("B11").Value = 0.0013
$ws.Range("B12").Value = 0.0086
$ws.Range("B13").Value = 0.0695
$ws.Range("B14").Value = 0.0008
$ws.Range("B15").Value = -1.2648
$ws.Range("B16").Value = -0.073
$ws.Range("B17").Value = 0.0006
$ws.Range("B18").Value = -0.0029
$ws.Range("B19").Value = 0.0186
$ws.Range("B20").Value = 0.542
$ws.Range("B21").Value = 0.031
$ws.Range("B22").Value = -0.0005
$ws.Range("B23").Value = 0.1513
$ws.Range("B24").Value = 0.1722
